$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (year headers): add 2019 in P4, 2020 in Q4, matching the
#     formatting used by the existing year header cells (D4:O4). ---
$ws.Range("N4").Copy()
$ws.Range("P4:Q4").PasteSpecial(-4122)
$ws.Range("P4").Value = 2019
$ws.Range("Q4").Value = 2020

# --- Row 5 (share of renewables, %): add 2019 value 35.67 in P5, and
#     an empty cell Q5 (2020 not yet available), using the same format
#     as the other "rounded" percentage cells (E5 / H5). ---
$ws.Range("E5").Copy()
$ws.Range("P5:Q5").PasteSpecial(-4122)
$ws.Range("P5").Value = 35.67

# --- Row 6 (hydropower output, mln kWh): add 2019 / 2020 values,
#     matching the formatting of the neighbouring O6 cell. ---
$ws.Range("O6").Copy()
$ws.Range("P6:Q6").PasteSpecial(-4122)
$ws.Range("P6").Value = 13859.3
$ws.Range("Q6").Value = 13979.1

# --- Update the active selection shown when the sheet is reopened. ---
$null = $ws.Range("P9").Select()
